# Commit: "Framework with active run flag includes: login module and Admin - user managemet"
# The test data row for TestCaseID=2 (row 3) has its RunFlag (column B) turned on,
# enabling that test case to run, and the active cell/selection moves to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1
$ws.Range("B3").Select()
